$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure price-column cells keep text (string) storage, matching the
# original inline-string cells, instead of Excel auto-coercing plain
# numeric-looking text into numbers.
$priceCells = @('D2', 'D3', 'D5', 'D8', 'D11', 'D12', 'D13', 'D15', 'D16', 'D17', 'D18', 'D19', 'D22', 'D25', 'D27', 'D34', 'D39', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50')
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '28.345.44'
$ws.Range('D3').Value = '1.582.52'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  -0.95%  '
$ws.Range('D5').Value = '213.66'
$ws.Range('E5').Value = '  +0.90%  '
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('E7').Value = '  -0.90%  '
$ws.Range('D8').Value = '23.72'
$ws.Range('E8').Value = '  +7.30%  '
$ws.Range('E9').Value = '  +0.94%  '
$ws.Range('E10').Value = '  -0.32%  '
$ws.Range('D11').Value = '0.0888'
$ws.Range('E11').Value = '  +2.11%  '
$ws.Range('D12').Value = '1.806.93'
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = '1.592.78'
$ws.Range('E13').Value = '  +1.81%  '
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = '0.529'
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('D16').Value = '28.324.74'
$ws.Range('E16').Value = '  +4.17%  '
$ws.Range('D17').Value = '63.86'
$ws.Range('E17').Value = '  +2.44%  '
$ws.Range('D18').Value = '232.46'
$ws.Range('E18').Value = '  +7.24%  '
$ws.Range('D19').Value = '0.0₃0708'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').Value = '4.12'
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('D25').Value = '151.56'
$ws.Range('E25').Value = '  -1.65%  '
$ws.Range('E26').Value = '  +1.04%  '
$ws.Range('D27').Value = '6.61'
$ws.Range('E27').Value = '  -1.18%  '
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('E29').Value = '  -0.85%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('E32').Value = '  -0.54%  '
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('D34').Value = '1.417.20'
$ws.Range('E34').Value = '  -2.68%  '
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('E36').Value = '  -5.62%  '
$ws.Range('E37').Value = '  -1.52%  '
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').Value = '2.53'
$ws.Range('E39').Value = '  +7.41%  '
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('E42').Value = '  -0.95%  '
$ws.Range('E43').Value = '  -2.40%  '
$ws.Range('D44').Value = '1.82'
$ws.Range('E44').Value = '  +5.12%  '
$ws.Range('D45').Value = '0.971'
$ws.Range('E45').Value = '  -3.03%  '
$ws.Range('D46').Value = '64.30'
$ws.Range('E46').Value = '  -0.56%  '
$ws.Range('D47').Value = '1.717.74'
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('D48').Value = '87.20'
$ws.Range('E48').Value = '  +1.55%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0526'
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0102'
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('E51').Value = '  +15.64%  '
